$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO" (was "Sheet1")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column before column C; this shifts the existing column C
# (header + data) to column D, carrying its formatting along.
$ws1.Columns.Item(3).Insert()

# New column C header + values: "Order Week" = old column A date + 6 days,
# formatted the same way column A's dates are.
$ws1.Cells.Item(1, 3).Value = "Order Week"
$ws1.Cells.Item(1, 1).Copy()
$ws1.Cells.Item(1, 3).PasteSpecial(-4122)

for ($r = 2; $r -le 14; $r++) {
    $origDate = $ws1.Cells.Item($r, 1).Value2

    # New column C ("Order Week") keeps the *original* ds date.
    $ws1.Cells.Item($r, 3).Value = $origDate
    $ws1.Cells.Item($r, 1).Copy()
    $ws1.Cells.Item($r, 3).PasteSpecial(-4122)

    # Column A is shifted forward 6 days (the new "ds").
    $ws1.Cells.Item($r, 1).Value = $origDate + 6

    # Column D (old "PO_Requested_Qty") is zeroed out - those quantities moved
    # to the new "Weekly Growth" sheet.
    $ws1.Cells.Item($r, 4).Value = 0
}

$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$ws2.Cells.Item(1, 1).Value = "ds"
$ws2.Cells.Item(1, 2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1, 3).Value = "Growth%"
$ws1.Cells.Item(1, 1).Copy()
$ws2.Cells.Item(1, 1).PasteSpecial(-4122)
$ws2.Cells.Item(1, 2).PasteSpecial(-4122)
$ws2.Cells.Item(1, 3).PasteSpecial(-4122)

$ws2.Cells.Item(2, 1).Value = 45572
$ws2.Cells.Item(2, 2).Value = 240
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(3, 1).Value = 45593
$ws2.Cells.Item(3, 2).Value = 108
$ws2.Cells.Item(3, 3).Value = -55.00000000000001

$ws1.Cells.Item(2, 1).Copy()
$ws2.Cells.Item(2, 1).PasteSpecial(-4122)
$ws2.Cells.Item(3, 1).PasteSpecial(-4122)

$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

$ws3.Cells.Item(1, 1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1, 2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1, 3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1, 4).Value = "Min_PO_Quantity"
$ws1.Cells.Item(1, 1).Copy()
$ws3.Cells.Item(1, 1).PasteSpecial(-4122)
$ws3.Cells.Item(1, 2).PasteSpecial(-4122)
$ws3.Cells.Item(1, 3).PasteSpecial(-4122)
$ws3.Cells.Item(1, 4).PasteSpecial(-4122)

$ws3.Cells.Item(2, 1).Value = 348
$ws3.Cells.Item(2, 2).Value = 174
$ws3.Cells.Item(2, 3).Value = 240
$ws3.Cells.Item(2, 4).Value = 108

$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"
$ws4.Outline.SummaryRow = 1
$ws4.Outline.SummaryColumn = 1
$ws4.PageSetup.LeftMargin = 54
$ws4.PageSetup.RightMargin = 54
$ws4.PageSetup.TopMargin = 72
$ws4.PageSetup.BottomMargin = 72
$ws4.PageSetup.HeaderMargin = 36
$ws4.PageSetup.FooterMargin = 36

$ws4.Cells.Item(1, 1).Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Cells.Item(1, 1).Copy()
$ws4.Cells.Item(1, 1).PasteSpecial(-4122)

$ws4.Cells.Item(2, 1).Value = 0

$ws1.Application.CutCopyMode = $false

# Put the focus back on the first sheet, matching the original activeTab.
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
